$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set for the "Negative Comments" mock data table.
# Row 1 is the header (unchanged); rows 2-13 hold the updated survey rows.
$data = @(
    @("CLINIC", "RESPONSE", "COMMENTS"),
    @("Radiology", "Unlikely", "Felt as if i was not a priority"),
    @("Special Care Baby Unit", "Unlikely", "Waited too long to find a parking spot"),
    @("Labour and Delivery Suite", "Extremely Unlikely", "Long wait times"),
    @("Labour and Delivery Suite", "Unlikely", "Long wait times"),
    @("Labour and Delivery Suite", "Extremely Unlikely", "Food was terrible"),
    @("Heart Failure", "Unlikely", "Clenliness isn't the best but otherwise okay"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "Clenliness isn't the best but otherwise okay"),
    @("Radiology", "Extremely Unlikely", "staff was rude"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "staff was rude"),
    @("Day Surgery", "Unlikely", "Waited too long to find a parking spot"),
    @("Gynaecology", "Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Day Surgery", "Extremely Unlikely", "staff tried to deal with me quickly rather than correctly. Not appropriate and i shouldve have been taken care of better. Would not recommend.")
)

# First clear out the old table (it previously spanned down to row 15).
$ws.Range("A1:C15").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Mirror the new selection left behind in the workbook (row after the table).
$ws.Range("A14:XFD102").Select() | Out-Null
